# Applies the per-row Coin/Link/Price/Volume(1h) updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.466.76"
$ws.Range("E2").Value = "  -5.92%  "
$ws.Range("D3").Value = "3.208.57"
$ws.Range("E3").Value = "  -8.92%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'576.02"
$ws.Range("E5").Value = "  -6.21%  "
$ws.Range("D6").Value = "'149.41"
$ws.Range("E6").Value = "  -13.86%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.202.16"
$ws.Range("E8").Value = "  -8.92%  "
$ws.Range("D9").Value = "'0.540"
$ws.Range("E9").Value = "  -11.42%  "
$ws.Range("D10").Value = "'0.170"
$ws.Range("E10").Value = "  -13.65%  "
$ws.Range("D11").Value = "'6.42"
$ws.Range("E11").Value = "  -11.34%  "
$ws.Range("D12").Value = "'0.495"
$ws.Range("E12").Value = "  -15.89%  "
$ws.Range("D13").Value = "'38.19"
$ws.Range("E13").Value = "  -17.85%  "
$ws.Range("D14").Value = "'0.0000240"
$ws.Range("E14").Value = "  -13.13%  "
$ws.Range("D15").Value = "3.711.65"
$ws.Range("E15").Value = "  -9.35%  "
$ws.Range("D16").Value = "66.403.89"
$ws.Range("E16").Value = "  -6.09%  "
$ws.Range("D17").Value = "3.200.86"
$ws.Range("E17").Value = "  -9.29%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'532.28"
$ws.Range("E18").Value = "  -13.81%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.114"
$ws.Range("E19").Value = "  -6.44%  "
$ws.Range("D20").Value = "'7.10"
$ws.Range("E20").Value = "  -15.86%  "
$ws.Range("D21").Value = "'14.97"
$ws.Range("E21").Value = "  -15.53%  "
$ws.Range("D22").Value = "'0.754"
$ws.Range("E22").Value = "  -14.67%  "
$ws.Range("D23").Value = "'7.68"
$ws.Range("E23").Value = "  -14.50%  "
$ws.Range("D24").Value = "'85.03"
$ws.Range("E24").Value = "  -13.69%  "
$ws.Range("D25").Value = "'13.20"
$ws.Range("E25").Value = "  -16.12%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'3.11"
$ws.Range("E27").Value = "  -17.95%  "
$ws.Range("D28").Value = "'2.15"
$ws.Range("E28").Value = "  -17.57%  "
$ws.Range("D29").Value = "'7.96"
$ws.Range("E29").Value = "  -13.12%  "
$ws.Range("D30").Value = "'28.90"
$ws.Range("E30").Value = "  -14.55%  "
$ws.Range("E31").Value = "  -16.84%  "
$ws.Range("D32").Value = "'1.12"
$ws.Range("E32").Value = "  -14.52%  "
$ws.Range("D33").Value = "'530.76"
$ws.Range("E33").Value = "  -14.54%  "
$ws.Range("D34").Value = "'6.47"
$ws.Range("E34").Value = "  -20.56%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'5.60"
$ws.Range("E35").Value = "  -18.08%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "'52.80"
$ws.Range("E37").Value = "  -7.43%  "
$ws.Range("D38").Value = "'0.0852"
$ws.Range("E38").Value = "  -15.08%  "
$ws.Range("D39").Value = "'0.0416"
$ws.Range("E39").Value = "  -13.91%  "
$ws.Range("D40").Value = "'9.02"
$ws.Range("E40").Value = "  -16.78%  "
$ws.Range("D41").Value = "'0.124"
$ws.Range("E41").Value = "  -14.50%  "
$ws.Range("D42").Value = "2.888.35"
$ws.Range("E42").Value = "  -14.33%  "
$ws.Range("D43").Value = "'2.60"
$ws.Range("E43").Value = "  -25.51%  "
$ws.Range("D44").Value = "0.0₃0582"
$ws.Range("E44").Value = "  -21.06%  "
$ws.Range("D45").Value = "'0.258"
$ws.Range("E45").Value = "  -17.32%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "'2.32"
$ws.Range("E47").Value = "  -21.07%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'25.69"
$ws.Range("E48").Value = "  -20.33%  "
$ws.Range("D49").Value = "'2.09"
$ws.Range("E49").Value = "  -18.68%  "
$ws.Range("D50").Value = "'0.113"
$ws.Range("E50").Value = "  -13.86%  "
$ws.Range("D51").Value = "'122.44"
$ws.Range("E51").Value = "  -8.60%  "

Write-Host "Applied 110 cell updates"
